$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 corresponds to student #1 (cccd=123456, hodem=Hoang).
# Column D is "ten" (first name) - rename it from "Nam" to "Giang".
$ws.Range("D2").Value = "Giang"
